$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global_Variables")

# Update the paystub options label and JSON value
$ws.Range("G1").Value = "paystub_A_options"
$ws.Range("G2").Value = '{"Rate" : 20 , "4_Digit_Account_Number" : 8698, "Numbe of Paystubs" : 5, "Period" : "Apr 01 2022"}'

# Move the active selection from G4 to F4
$ws.Range("F4").Select()
